$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.5316481590271
$ws.Range("B1").Value = 2.866567373275757
$ws.Range("C1").Value = 2.071220636367798
$ws.Range("D1").Value = 1.913193345069885
$ws.Range("E1").Value = 1.974348902702332
